# Updated test data for DC, TripCurrent, Voltdrop, BatteryStandby
$wb = $excel.ActiveWorkbook

# --- GreenColorPercentage sheet ---
$wsGreen = $wb.Worksheets.Item("GreenColorPercentage")
$wsGreen.Activate()
$wsGreen.Range("B8").Value = 30
$wsGreen.Range("B4").Select()

# --- VDWorstCaseYellowPercentage sheet ---
$wsVD = $wb.Worksheets.Item("VDWorstCaseYellowPercentage")
$wsVD.Activate()
$wsVD.Range("B8").Value = 31
$wsVD.Range("B6").Select()

# --- VtgDropYellowColorPercentage sheet ---
$wsVtg = $wb.Worksheets.Item("VtgDropYellowColorPercentage")
$wsVtg.Activate()
$wsVtg.Range("B8").Value = 51
$wsVtg.Range("B8").Select()

# --- RedColorPercentage sheet ---
$wsRed = $wb.Worksheets.Item("RedColorPercentage")
$wsRed.Activate()
$wsRed.Range("B8").Value = 52
$wsRed.Range("B4").Select()

# Return to the GreenColorPercentage sheet, which ends up the active/selected tab
$wsGreen.Activate()
$wsGreen.Range("B4").Select()
